# Fruta / hortaliza, semanal
# A new weekly price record is inserted at row 35, pushing the existing
# rows 35-76 down to 36-77 (dimension grows from A1:R76 to A1:R77).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 35 (shifts 35..76 -> 36..77)
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new weekly entry
$ws.Range("A35").Value = 5
$ws.Range("B35").Value = "Macroferia Regional de Talca"
$ws.Range("C35").Value = "Maule"
$ws.Range("D35").Value = 44799
$ws.Range("E35").Value = 7
$ws.Range("F35").Value = 100112026
$ws.Range("G35").Value = "Haba"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 150
$ws.Range("K35").Value = 12000
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = 12000
$ws.Range("N35").Value = "$/saco 25 kilos"
$ws.Range("O35").Value = "Provincia del Elquí"
$ws.Range("P35").Value = 480
$ws.Range("Q35").Value = 25
$ws.Range("R35").Value = "Hortaliza"
